$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to Text format first so the numeric-looking
# strings are preserved exactly (matching the source data's text values)
# instead of being auto-converted to floating point numbers.
# (Each cell is set individually since a single comma-joined multi-area
# Range() reference only reliably formats its first area.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "40.812.42"
$ws.Range("E2").Value = "  +3.71%  "

$ws.Range("D3").Value = "2.215.25"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "229.53"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("E6").Value = "  +1.89%  "

$ws.Range("E7").Value = "  +0.66%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.404"
$ws.Range("E9").Value = "  +1.76%  "

$ws.Range("D10").Value = "0.0867"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("D12").Value = "2.541.52"
$ws.Range("E12").Value = "  +2.49%  "

$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "0.821"
$ws.Range("E15").Value = "  +0.98%  "

$ws.Range("D16").Value = "5.61"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("D17").Value = "2.208.42"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").Value = "40.679.94"
$ws.Range("E18").Value = "  +3.54%  "

$ws.Range("D19").Value = "74.01"
$ws.Range("E19").Value = "  +3.01%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.17"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +5.88%  "

$ws.Range("D22").Value = "250.33"
$ws.Range("E22").Value = "  +8.38%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -8.83%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "173.07"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  +1.49%  "

$ws.Range("D28").Value = "0.144"
$ws.Range("E28").Value = "  +3.78%  "

$ws.Range("D29").Value = "20.31"
$ws.Range("E29").Value = "  +1.94%  "

$ws.Range("E30").Value = "  +2.29%  "

$ws.Range("D31").Value = "2.83"
$ws.Range("E31").Value = "  +5.20%  "

$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("E36").Value = "  +2.33%  "

$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +6.55%  "

$ws.Range("E38").Value = "  +1.79%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").Value = "4.94"
$ws.Range("E40").Value = "  +14.55%  "

$ws.Range("E41").Value = "  +1.73%  "

$ws.Range("E42").Value = "  +9.18%  "

$ws.Range("D43").Value = "101.34"
$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("E44").Value = "  +4.70%  "

$ws.Range("D45").Value = "17.39"
$ws.Range("E45").Value = "  -2.41%  "

$ws.Range("D46").Value = "1.511.58"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("D47").Value = "0.0940"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("D50").Value = "0.000206"
$ws.Range("E50").Value = "  +38.45%  "

$ws.Range("D51").Value = "9.62"
$ws.Range("E51").Value = "  +11.55%  "

